$wb = $excel.ActiveWorkbook

# --- Budget sheet: fiscal year bump 2022 -> 2023 ---
$wsBudget = $wb.Worksheets.Item("Budget")
$wsBudget.Range("A2").Value = 2023
[void]$wsBudget.Range("I14").Select()

# --- Criteria sheet: rename INTERNETREPORT field to INTERNET_REPORT ---
$wsCriteria = $wb.Worksheets.Item("Criteria")
$wsCriteria.Range("B2").Value = "[INTERSTATE]=|Y| AND [INTERNET_REPORT]=|State|"
$wsCriteria.Range("B3").Value = "[INTERSTATE]='Y' AND [INTERNET_REPORT]='State'"

# Style on the criteria text cells was simplified back to the default/Normal style
$wsCriteria.Range("B2:B3").Style = "Normal"
[void]$wsCriteria.Range("C3").Select()

# Restore Budget as the active/selected sheet and cell
[void]$wsBudget.Select()
[void]$wsBudget.Range("I14").Select()
